$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refreshed prices and 1h volume % changes; two coin rows (Polkadot/BitcoinCash,
# ImmutableX/EthereumClassic, OKB/Bittensor, TheGraph/Monero) swapped order; ThetaToken -> BitgetToken
$ws.Range("D2").Value = "64.011.79"
$ws.Range("E2").Value = "  -4.31%  "
$ws.Range("D3").Value = "3.064.98"
$ws.Range("E3").Value = "  -4.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.81"
$ws.Range("E5").Value = "  -8.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.19"
$ws.Range("E6").Value = "  -12.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "3.060.87"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -8.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -11.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.80"
$ws.Range("E11").Value = "  -8.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -7.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  -10.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.98"
$ws.Range("E14").Value = "  -15.01%  "
$ws.Range("D15").Value = "3.576.60"
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("D16").Value = "64.264.60"
$ws.Range("E16").Value = "  -3.85%  "
$ws.Range("D17").Value = "3.111.39"
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.111"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "481.78"
$ws.Range("E19").Value = "  -9.17%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.40"
$ws.Range("E20").Value = "  -8.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.39"
$ws.Range("E21").Value = "  -7.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.667"
$ws.Range("E22").Value = "  -11.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.93"
$ws.Range("E23").Value = "  -11.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.18"
$ws.Range("E24").Value = "  -9.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.94"
$ws.Range("E25").Value = "  -11.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.986"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("E27").Value = "  -10.27%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.78"
$ws.Range("E28").Value = "  -10.43%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.89"
$ws.Range("E29").Value = "  -12.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("E30").Value = "  -12.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.35"
$ws.Range("E32").Value = "  -10.71%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "53.57"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "464.55"
$ws.Range("E35").Value = "  -12.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.63"
$ws.Range("E36").Value = "  -12.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.87"
$ws.Range("E37").Value = "  -14.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0387"
$ws.Range("E38").Value = "  -8.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0759"
$ws.Range("E39").Value = "  -11.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.04"
$ws.Range("E40").Value = "  -12.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  -9.55%  "
$ws.Range("D42").Value = "2.708.05"
$ws.Range("E42").Value = "  -6.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("E44").Value = "  -17.85%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "120.24"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.229"
$ws.Range("E46").Value = "  -12.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.103"
$ws.Range("E47").Value = "  -9.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.71"
$ws.Range("E48").Value = "  -12.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.85"
$ws.Range("E49").Value = "  -13.52%  "
$ws.Range("D50").Value = "0.0₃0474"
$ws.Range("E50").Value = "  -18.87%  "
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.22"
$ws.Range("E51").Value = "  -1.81%  "
